$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (Leve Item ID=5471) - You Grow, Girl | Growth Formula Beta
$ws.Range("H18").Value = 1825.25
$ws.Range("I18").Value = 1825.25
$ws.Range("K18").Value = 1825.25
$ws.Range("M18").Value = -1541.25

# Row 40 (Leve Item ID=5505) - Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 1999.9
$ws.Range("J40").Value = 1999.5
$ws.Range("L40").Value = 1999.5
$ws.Range("N40").Value = -2349.5

# Row 55 (Leve Item ID=5517) - A Real Smooth Move | Lanolin
$ws.Range("H55").Value = 760.44446
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 863.4286
$ws.Range("K55").Value = 400
$ws.Range("L55").Value = 863.4286
$ws.Range("M55").Value = -186
$ws.Range("N55").Value = -1291.4286

# Row 99 (Leve Item ID=19883) - Rumor Has It | Commanding Craftsman's Tea
$ws.Range("H99").Value = 197.75
$ws.Range("I99").Value = 197.75
$ws.Range("K99").Value = 593.25
$ws.Range("M99").Value = 904.75

# Row 132 (Leve Item ID=44049) - Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 2106.8125
$ws.Range("I132").Value = 2111.6785
$ws.Range("J132").Value = 2072.75
$ws.Range("K132").Value = 6335.0355
$ws.Range("L132").Value = 6218.25
$ws.Range("M132").Value = -3805.0355
$ws.Range("N132").Value = -11278.25

# Row 137 (Leve Item ID=44013) - Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 4619.0557
$ws.Range("I137").Value = 1725.75
$ws.Range("K137").Value = 5177.25
$ws.Range("M137").Value = -2627.25

# Row 138 (Leve Item ID=44169) - All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 6898.44
$ws.Range("I138").Value = 7962.2144
$ws.Range("K138").Value = 23886.6432
$ws.Range("M138").Value = -18746.6432

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID=44147) - Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 15025.246
$ws.Range("I32").Value = 6213.5117
$ws.Range("J32").Value = 27655.4
$ws.Range("K32").Value = 6213.5117
$ws.Range("L32").Value = 27655.4
$ws.Range("M32").Value = -5926.5117
$ws.Range("N32").Value = -28229.4

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID=19939) - High Steal | High Steel Nugget
$ws.Range("H94").Value = 724.8570999999999
$ws.Range("I94").Value = 654.8
$ws.Range("K94").Value = 654.8
$ws.Range("M94").Value = -203.8

$ws = $wb.Worksheets.Item("CRP")
# Row 41 (Leve Item ID=1917) - The Lone Bowman | Oak Longbow
$ws.Range("H41").Value = 19250
$ws.Range("I41").Value = 10250
$ws.Range("J41").Value = 20750
$ws.Range("K41").Value = 10250
$ws.Range("L41").Value = 20750
$ws.Range("M41").Value = -9822
$ws.Range("N41").Value = -21606

# Row 62 (Leve Item ID=12580) - Splinter in the Sewers | Cedar Lumber
$ws.Range("H62").Value = 54365.375
$ws.Range("I62").Value = 4984.6665
$ws.Range("K62").Value = 4984.6665
$ws.Range("M62").Value = -4360.6665

# Row 65 (Leve Item ID=12580) - The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Range("H65").Value = 54365.375
$ws.Range("I65").Value = 4984.6665
$ws.Range("K65").Value = 24923.3325
$ws.Range("M65").Value = -21803.3325

# Row 68 (Leve Item ID=10611) - Do You Even String Bow | Holy Cedar Composite Bow
$ws.Range("H68").Value = 53499.5
$ws.Range("J68").Value = 53499.5
$ws.Range("L68").Value = 53499.5
$ws.Range("N68").Value = -54997.5

# Row 71 (Leve Item ID=10611) - Win One Bow, Get Three Free (L) | Holy Cedar Composite Bow
$ws.Range("H71").Value = 53499.5
$ws.Range("J71").Value = 53499.5
$ws.Range("L71").Value = 160498.5
$ws.Range("N71").Value = -167986.5

# Row 74 (Leve Item ID=10636) - License to Heal | Dark Chestnut Rod
$ws.Range("H74").Value = 77721.5
$ws.Range("J74").Value = 77721.5
$ws.Range("L74").Value = 77721.5
$ws.Range("N74").Value = -79469.5

# Row 77 (Leve Item ID=10636) - Purified Polyrhythm (L) | Dark Chestnut Rod
$ws.Range("H77").Value = 77721.5
$ws.Range("J77").Value = 77721.5
$ws.Range("L77").Value = 233164.5
$ws.Range("N77").Value = -241900.5

# Row 107 (Leve Item ID=27689) - Built to Last | White Oak Lumber
$ws.Range("H107").Value = 620.6667
$ws.Range("I107").Value = 368.75
$ws.Range("K107").Value = 368.75
$ws.Range("M107").Value = 1551.25

# Row 132 (Leve Item ID=44019) - Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 812.2105
$ws.Range("I132").Value = 812.2105
$ws.Range("K132").Value = 2436.6315
$ws.Range("M132").Value = 93.36850000000004

$ws = $wb.Worksheets.Item("CUL")
# Row 14 (Leve Item ID=12886) - Keep Your Powder Dry | Kukuru Powder
$ws.Range("H14").Value = 2477.375
$ws.Range("I14").Value = 2477.375
$ws.Range("K14").Value = 7432.125
$ws.Range("M14").Value = -7259.125

# Row 33 (Leve Item ID=4867) - Cooking with Gas | Chicken Stock
$ws.Range("H33").Value = 849.2308
$ws.Range("J33").Value = 104.42857
$ws.Range("L33").Value = 626.57142
$ws.Range("N33").Value = -1192.57142

# Row 54 (Leve Item ID=4671) - Good Eats in Ishgard | Salt Cod Puffs
$ws.Range("H54").Value = 3000
$ws.Range("J54").Value = 3000
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10118

# Row 68 (Leve Item ID=12895) - Such a Butter Face | Fermented Butter
$ws.Range("H68").Value = 2986.1
$ws.Range("J68").Value = 2968.2856
$ws.Range("L68").Value = 8904.856800000001
$ws.Range("N68").Value = -10526.8568

# Row 71 (Leve Item ID=12895) - No Margarine of Error (L) | Fermented Butter
$ws.Range("H71").Value = 2986.1
$ws.Range("J71").Value = 2968.2856
$ws.Range("L71").Value = 26714.5704
$ws.Range("N71").Value = -34826.5704

# Row 74 (Leve Item ID=12859) - The Nutcracker's Sweets | Royal Eggs
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = 0

# Row 77 (Leve Item ID=12859) - Time for a Midnight Snack (L) | Royal Eggs
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = 0

# Row 80 (Leve Item ID=12890) - Saucy for a Suitor | Hollandaise Sauce
$ws.Range("H80").Value = 5955.9
$ws.Range("I80").Value = 5926.5
$ws.Range("K80").Value = 17779.5
$ws.Range("M80").Value = -16843.5

# Row 83 (Leve Item ID=12890) - Saved by the Sauce (L) | Hollandaise Sauce
$ws.Range("H83").Value = 5955.9
$ws.Range("I83").Value = 5926.5
$ws.Range("K83").Value = 53338.5
$ws.Range("M83").Value = -48658.5

# Row 113 (Leve Item ID=27843) - Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 2387.3
$ws.Range("J113").Value = 2167.8572
$ws.Range("L113").Value = 6503.571599999999
$ws.Range("N113").Value = -10843.5716

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID=5062) - Copper and Robbers | Copper Ingot
$ws.Range("H2").Value = 500.46667
$ws.Range("J2").Value = 1519.5
$ws.Range("L2").Value = 1519.5
$ws.Range("N2").Value = -1745.5

# Row 33 (Leve Item ID=4450) - Thaumaturge Is Magic | Fluorite Ring
$ws.Range("H33").Value = 15996.5
$ws.Range("J33").Value = 15996.5
$ws.Range("L33").Value = 15996.5
$ws.Range("N33").Value = -16500.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID=5277) - Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 1262.3846
$ws.Range("I22").Value = 744.4
$ws.Range("K22").Value = 744.4
$ws.Range("M22").Value = -449.4

# Row 27 (Leve Item ID=5277) - Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 1262.3846
$ws.Range("I27").Value = 744.4
$ws.Range("K27").Value = 744.4
$ws.Range("M27").Value = -637.4

# Row 40 (Leve Item ID=36248) - Best Served Toad | Toad Leather
$ws.Range("H40").Value = 5465
$ws.Range("I40").Value = 5558.2
$ws.Range("J40").Value = 4999
$ws.Range("K40").Value = 5558.2
$ws.Range("L40").Value = 4999
$ws.Range("M40").Value = -5422.2
$ws.Range("N40").Value = -5271

# Row 46 (Leve Item ID=5282) - Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 2996.4285
$ws.Range("I46").Value = 1993.75
$ws.Range("J46").Value = 4333.3335
$ws.Range("K46").Value = 1993.75
$ws.Range("L46").Value = 4333.3335
$ws.Range("M46").Value = -1805.75
$ws.Range("N46").Value = -4709.3335

# Row 55 (Leve Item ID=5284) - It's Not a Job, It's a Calling | Peiste Leather
$ws.Range("H55").Value = 707.3333
$ws.Range("I55").Value = 294
$ws.Range("K55").Value = 294
$ws.Range("M55").Value = -121

# Row 100 (Leve Item ID=19995) - Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 5085.143
$ws.Range("I100").Value = 1749
$ws.Range("K100").Value = 1749
$ws.Range("M100").Value = -1208

# Row 122 (Leve Item ID=36247) - Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 4316.3335
$ws.Range("I122").Value = 3974
$ws.Range("K122").Value = 11922
$ws.Range("M122").Value = -9472

$ws = $wb.Worksheets.Item("WVR")
# Row 86 (Leve Item ID=11977) - Felt for the Fallen | Chimerical Felt
$ws.Range("H86").Value = 40000
$ws.Range("J86").Value = 40000
$ws.Range("L86").Value = 40000
$ws.Range("N86").Value = -42246

# Row 89 (Leve Item ID=11977) - Blinded Veil of Vigilance (L) | Chimerical Felt
$ws.Range("H89").Value = 40000
$ws.Range("J89").Value = 40000
$ws.Range("L89").Value = 200000
$ws.Range("N89").Value = -211232

# Row 113 (Leve Item ID=27752) - A Tender Table | Pixie Floss
$ws.Range("H113").Value = 1301.0667
$ws.Range("I113").Value = 801.3333
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 2403.9999
$ws.Range("L113").Value = 9900
$ws.Range("M113").Value = -233.9998999999998
$ws.Range("N113").Value = -14240

# Row 122 (Leve Item ID=36208) - Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 6320.6
$ws.Range("I122").Value = 2400.875
$ws.Range("K122").Value = 7202.625
$ws.Range("M122").Value = -4752.625

# Row 132 (Leve Item ID=44029) - Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2533.4
$ws.Range("I132").Value = 1859.6154
$ws.Range("J132").Value = 3784.7144
$ws.Range("K132").Value = 5578.8462
$ws.Range("L132").Value = 11354.1432
$ws.Range("M132").Value = -3048.8462
$ws.Range("N132").Value = -16414.1432
